$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.653.92"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.653.63"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9966"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3626"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.96"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3262"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.135"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07059"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.020"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.56"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.659.37"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.621"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -7.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06615"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9979"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.90"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.917"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -8.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.59"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.712.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.434"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.391"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -13.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.66"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.77%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.841.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.44%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.211"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.96"
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.861"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -12.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08434"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.661"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.29"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.66%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.283"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.216"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.20%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06027"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.98%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02232"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2071"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.186"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9975"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5925"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.821"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5639"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.945"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06975"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.190"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.44%  "
